$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("survey")
$ws2 = $wb.Worksheets.Item("choices")

# "select_one yes_no" -> "boolean" for the three relevant survey rows
$ws1.Range("A4").Value = "boolean"
$ws1.Range("A5").Value = "boolean"
$ws1.Range("A6").Value = "boolean"

# the yes_no choice list is no longer needed - remove it from the choices sheet
$ws2.Rows("3:4").Delete()

# restore the selections recorded in the edited workbook
$ws2.Activate()
[void]$ws2.Range("D14").Select()
$ws1.Activate()
[void]$ws1.Range("C21").Select()
